$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.081.16'
$ws.Range("E2").Value = '  -0.14%  '

# Row 3
$ws.Range("D3").Value = '1.759.00'
$ws.Range("E3").Value = '  +2.40%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9968'
$ws.Range("E4").Value = '  -0.45%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.77'
$ws.Range("E5").Value = '  +1.35%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9949'
$ws.Range("E6").Value = '  -0.54%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5180'
$ws.Range("E7").Value = '  +9.99%  '

# Row 8
$ws.Range("E8").Value = '  +4.78%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.30'
$ws.Range("E9").Value = '  +0.30%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07276'
$ws.Range("E10").Value = '  +0.14%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.071'
$ws.Range("E11").Value = '  +2.83%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9929'
$ws.Range("E12").Value = '  -0.70%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.47'
$ws.Range("E13").Value = '  +2.82%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.012'
$ws.Range("E14").Value = '  +2.80%  '

# Row 15
$ws.Range("D15").Value = '1.747.92'
$ws.Range("E15").Value = '  +1.93%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.864'
$ws.Range("E16").Value = '  -0.34%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.37'
$ws.Range("E17").Value = '  -1.66%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001036'
$ws.Range("E18").Value = '  -0.10%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06409'
$ws.Range("E19").Value = '  +0.87%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9976'
$ws.Range("E20").Value = '  -0.23%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.67'
$ws.Range("E21").Value = '  +1.07%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.785'
$ws.Range("E22").Value = '  +3.04%  '

# Row 23
$ws.Range("D23").Value = '27.180.57'
$ws.Range("E23").Value = '  +0.02%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.29'
$ws.Range("E24").Value = '  +4.17%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.047'
$ws.Range("E25").Value = '  -4.57%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.27'
$ws.Range("E26").Value = '  -2.52%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.07'
$ws.Range("E27").Value = '  +3.12%  '

# Row 28
$ws.Range("D28").Value = '1.943.74'
$ws.Range("E28").Value = '  +1.67%  '

# Row 29
$ws.Range("E29").Value = '  +7.22%  '

# Row 30
$ws.Range("E30").Value = '  +0.48%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.050'
$ws.Range("E31").Value = '  +2.85%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09587'
$ws.Range("E32").Value = '  +4.51%  '

# Row 33
$ws.Range("E33").Value = '  +0.09%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.429'
$ws.Range("E34").Value = '  +2.07%  '

# Row 35
$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02199'
$ws.Range("E35").Value = '  +0.17%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05921'
$ws.Range("E36").Value = '  +1.80%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.07'
$ws.Range("E37").Value = '  +1.06%  '

# Row 38
$ws.Range("B38").Value = 'InternetComputer(DFINITY)'
$ws.Range("C38").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.795'
$ws.Range("E38").Value = '  +1.36%  '

# Row 39
$ws.Range("B39").Value = 'TheSandbox'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.6088'
$ws.Range("E39").Value = '  +3.30%  '

# Row 40
$ws.Range("B40").Value = 'Algorand'
$ws.Range("C40").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2004'
$ws.Range("E40").Value = '  +0.57%  '

# Row 41
$ws.Range("B41").Value = 'WEMIXTOKEN'
$ws.Range("C41").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.424'
$ws.Range("E41").Value = '  +2.26%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.114'
$ws.Range("E42").Value = '  -0.77%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.711'
$ws.Range("E43").Value = '  +3.70%  '

# Row 44
$ws.Range("E44").Value = '  +3.27%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.600'
$ws.Range("E45").Value = '  +1.40%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5690'
$ws.Range("E46").Value = '  +0.91%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '120.49'
$ws.Range("E47").Value = '  +2.29%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.861'
$ws.Range("E48").Value = '  +1.30%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06679'
$ws.Range("E49").Value = '  +0.46%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.102'
$ws.Range("E50").Value = '  +1.60%  '

# Row 51
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '69.97'
$ws.Range("E51").Value = '  +0.10%  '
